$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 176.6
$ws.Range("I6").Value = 176.6
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 529.8
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = -417.8
$ws.Range("M6").ClearContents()

$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("L16").ClearContents()

$ws.Range("H33").Value = 1430.3077
$ws.Range("I33").Value = 1590.9
$ws.Range("K33").Value = 1590.9
$ws.Range("M33").Value = -1361.9

$ws.Range("H86").Value = 41669590
$ws.Range("I86").Value = 53746310
$ws.Range("K86").Value = 53746310
$ws.Range("M86").Value = -53745187

$ws.Range("H89").Value = 41669590
$ws.Range("I89").Value = 53746310
$ws.Range("K89").Value = 268731550
$ws.Range("M89").Value = -268725934

$ws.Range("H132").Value = 1247.449
$ws.Range("I132").Value = 1191.1364
$ws.Range("K132").Value = 3573.4092
$ws.Range("M132").Value = -1043.4092

$ws.Range("H137").Value = 2410.8845
$ws.Range("I137").Value = 2442.818
$ws.Range("K137").Value = 7328.454000000001
$ws.Range("M137").Value = -4778.454000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 29804.334
$ws.Range("I28").Value = 17735.5
$ws.Range("K28").Value = 17735.5
$ws.Range("M28").Value = -17543.5

$ws.Range("H32").Value = 3719.77
$ws.Range("I32").Value = 3831.0833
$ws.Range("J32").Value = 1048.25
$ws.Range("K32").Value = 3831.0833
$ws.Range("L32").Value = 1048.25
$ws.Range("M32").Value = -3544.0833
$ws.Range("N32").Value = -1622.25

$ws.Range("H61").Value = 5271.491
$ws.Range("J61").Value = 13327.615
$ws.Range("L61").Value = 13327.615
$ws.Range("N61").Value = -13751.615

$ws.Range("H74").Value = 18668.885
$ws.Range("I74").Value = 22936.963
$ws.Range("J74").Value = 4264.125
$ws.Range("K74").Value = 22936.963
$ws.Range("L74").Value = 4264.125
$ws.Range("M74").Value = -22062.963
$ws.Range("N74").Value = -6012.125

$ws.Range("H77").Value = 18668.885
$ws.Range("I77").Value = 22936.963
$ws.Range("J77").Value = 4264.125
$ws.Range("K77").Value = 114684.815
$ws.Range("L77").Value = 21320.625
$ws.Range("M77").Value = -110316.815
$ws.Range("N77").Value = -30056.625

$ws.Range("H99").Value = 29804.334
$ws.Range("I99").Value = 17735.5
$ws.Range("K99").Value = 17735.5
$ws.Range("M99").Value = -14740.5

$ws.Range("H132").Value = 7188.8716
$ws.Range("I132").Value = 6038.8887
$ws.Range("K132").Value = 18116.6661
$ws.Range("M132").Value = -15586.6661

$ws.Range("H136").Value = 5271.491
$ws.Range("J136").Value = 13327.615
$ws.Range("L136").Value = 39982.845
$ws.Range("N136").Value = -45082.845

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 2550
$ws.Range("I11").Value = 100
$ws.Range("J11").Value = 5000
$ws.Range("K11").Value = 100
$ws.Range("L11").Value = 5000
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = -5280

$ws.Range("H22").Value = 340.66666
$ws.Range("I22").Value = 289
$ws.Range("J22").Value = 444
$ws.Range("K22").Value = 289
$ws.Range("L22").Value = 444
$ws.Range("M22").Value = -116
$ws.Range("N22").Value = -790

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 227.5
$ws.Range("I22").Value = 199.5
$ws.Range("J22").Value = 255.5
$ws.Range("K22").Value = 199.5
$ws.Range("L22").Value = 255.5
$ws.Range("M22").Value = 150.5
$ws.Range("N22").Value = -955.5

$ws.Range("H58").Value = 7940810
$ws.Range("I58").Value = 11365000
$ws.Range("K58").Value = 11365000
$ws.Range("M58").Value = -11364797

$ws.Range("H64").Value = 72489
$ws.Range("J64").Value = 72489
$ws.Range("L64").Value = 72489
$ws.Range("N64").Value = -72985

$ws.Range("H67").Value = 72489
$ws.Range("J67").Value = 72489
$ws.Range("L67").Value = 72489
$ws.Range("N67").Value = -74205

$ws.Range("H132").Value = 6534.3823
$ws.Range("J132").Value = 9614.944
$ws.Range("L132").Value = 28844.832
$ws.Range("N132").Value = -33904.83199999999

$ws.Range("H136").Value = 7940810
$ws.Range("I136").Value = 11365000
$ws.Range("K136").Value = 34095000
$ws.Range("M136").Value = -34092450

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1616.6666
$ws.Range("I131").Value = 1447.375
$ws.Range("J131").Value = 1810.1428
$ws.Range("K131").Value = 4342.125
$ws.Range("L131").Value = 5430.428400000001
$ws.Range("M131").Value = 697.875
$ws.Range("N131").Value = -15510.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H122").Value = 3412841.8
$ws.Range("I122").Value = 3980760
$ws.Range("J122").Value = 5332.3335
$ws.Range("K122").Value = 11942280
$ws.Range("L122").Value = 15997.0005
$ws.Range("M122").Value = -11939830
$ws.Range("N122").Value = -20897.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3076.6
$ws.Range("I22").Value = 136.5
$ws.Range("J22").Value = 5036.6665
$ws.Range("K22").Value = 136.5
$ws.Range("L22").Value = 5036.6665
$ws.Range("M22").Value = 158.5
$ws.Range("N22").Value = -5626.6665

$ws.Range("H27").Value = 3076.6
$ws.Range("I27").Value = 136.5
$ws.Range("J27").Value = 5036.6665
$ws.Range("K27").Value = 136.5
$ws.Range("L27").Value = 5036.6665
$ws.Range("M27").Value = -29.5
$ws.Range("N27").Value = -5250.6665

$ws.Range("H29").Value = 1550
$ws.Range("I29").Value = 1550
$ws.Range("K29").Value = 1550
$ws.Range("M29").Value = -1255

$ws.Range("H95").Value = 200000
$ws.Range("J95").Value = 200000
$ws.Range("L95").Value = 200000
$ws.Range("N95").Value = -205492

$ws.Range("H112").Value = 46187
$ws.Range("J112").Value = 46187
$ws.Range("L112").Value = 46187
$ws.Range("N112").Value = -49141

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 40001600
$ws.Range("I81").Value = 2000
$ws.Range("K81").Value = 4000
$ws.Range("M81").Value = -2939

$ws.Range("H84").Value = 40001600
$ws.Range("I84").Value = 2000
$ws.Range("K84").Value = 20000
$ws.Range("M84").Value = -14696
